$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "Delbert"
$ws.Range("B11").Value = "Parker"
$ws.Range("C11").Value = "ronald.kshlerin@yahoo.com"

# phoneNumber must stay text (matches the other rows, which store it as a
# shared string, not a number) -- force it with a leading apostrophe and
# then strip the resulting "Quote Prefix" cell format back to Normal so the
# cell ends up with the same (default) style as its neighbours.
$ws.Range("D11").Value = "'6056110055"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "Scientist"
$ws.Range("F11").Value = "4rZFBc65"
